$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels in row 1 (prediction / error columns get a "day count"
# wording change: "N 预测值" -> "(N-1) 天预测值", "N 误差" -> "(N-1) 天误差")
$ws.Range("E1").Value = "1 天预测值"
$ws.Range("F1").Value = "2 天预测值"
$ws.Range("G1").Value = "3 天预测值"
$ws.Range("H1").Value = "4 天预测值"
$ws.Range("I1").Value = "5 天预测值"
$ws.Range("J1").Value = "6 天预测值"
$ws.Range("K1").Value = "7 天预测值"

$ws.Range("M1").Value = "1 天误差"
$ws.Range("N1").Value = "2 天误差"
$ws.Range("O1").Value = "3 天误差"
$ws.Range("P1").Value = "4 天误差"
$ws.Range("Q1").Value = "5 天误差"
$ws.Range("R1").Value = "6 天误差"
$ws.Range("S1").Value = "7 天误差"

# Widen columns E:K so the new, longer labels are readable (~10.875 characters)
$ws.Range("E1:K1").EntireColumn.ColumnWidth = 10.14

# Make "Sheet1" the active/selected tab instead of "Chart1"
$ws.Activate()
